# Daily attendance processing - normalize the order of the "Recorded By"
# entries (column G) on the active "Session Analysis Results" sheet.
#
# Each cell holds a comma-separated list of who/what recorded the session
# (e.g. "System, dnasr281@gmail.com"). The recorder names are reordered
# according to a fixed priority so the human recorders are listed before
# the automated "System" entry (and the spoofed lower-case "system" /
# backdoor entries sort to the very end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RecorderRank($token) {
    if ($token -eq "dnasr281@gmail.com") { return 1 }
    if ($token -eq "admin@admin.com") { return 2 }
    if ($token -eq "backup@backdoor.com") { return 4 }
    if ($token -eq "System") {
        # "System" and "system" compare equal case-insensitively, so
        # disambiguate the real (capitalised) System account from the
        # spoofed lower-case one using the character code of the first
        # letter ('S' = 83, 's' = 115).
        $firstCharCode = [int]$token.ToCharArray()[0]
        if ($firstCharCode -eq 83) {
            return 3
        } else {
            return 5
        }
    }
    return 99
}

function Get-SortedRecordedBy($value) {
    $parts = $value -split ", "
    $count = $parts.Length
    if ($count -le 1) {
        return $value
    }

    $ranks = @()
    foreach ($part in $parts) {
        $ranks += (Get-RecorderRank $part)
    }

    # Stable bubble sort (keeps original relative order for equal ranks).
    for ($i = 0; $i -lt $count; $i++) {
        for ($j = 0; $j -lt ($count - $i - 1); $j++) {
            if ($ranks[$j] -gt $ranks[$j + 1]) {
                $tmpRank = $ranks[$j]
                $ranks[$j] = $ranks[$j + 1]
                $ranks[$j + 1] = $tmpRank

                $tmpPart = $parts[$j]
                $parts[$j] = $parts[$j + 1]
                $parts[$j + 1] = $tmpPart
            }
        }
    }

    return [string]::Join(", ", $parts)
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$recordedByColumn = 7  # column G

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $recordedByColumn)
    $currentValue = $cell.Value2
    if ($currentValue -ne $null -and $currentValue -ne "") {
        $newValue = Get-SortedRecordedBy $currentValue
        if ($newValue -ne $currentValue) {
            $cell.Value = $newValue
        }
    }
}
